# "[plugin-excel] Add ability to set collection of ranges"
#
# The test template workbook gets:
#  1. its named ranges used for worksheet data connections renamed
#     (each gets a numeric "1" suffix appended, e.g. "...Node_Media" ->
#     "...Node_Media1"), reflecting the new collection-of-ranges support;
#  2. a new data row appended to the "RepeatingData" sheet (name1/status1);
#  3. "RepeatingData" becomes the active sheet/selection instead of
#     "DifferentTypes".

$wb = $excel.ActiveWorkbook

# 1. Rename every defined name (the hidden "_xlcn.WorksheetConnection_..."
#    names backing the worksheet data connections) by appending "1".
foreach ($n in $wb.Names) {
    $n.Name = $n.Name + "1"
}

# 2. Append a new row of sample data to the RepeatingData sheet.
$ws3 = $wb.Worksheets.Item("RepeatingData")
$ws3.Range("A9").Value = "name1"
$ws3.Range("B9").Value = "status1"

# 3. Make RepeatingData the active sheet and select A5 on it (previously
#    DifferentTypes was the active/selected sheet).
$ws3.Activate() | Out-Null
$ws3.Range("A5").Select() | Out-Null
